$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new employee-import column layout ---------------
$ws.Range("A1").Value = "comp_code"
$ws.Range("B1").Value = "branch_code"
$ws.Range("C1").Value = "nik"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "first_name"
$ws.Range("F1").Value = "last_name"
$ws.Range("G1").Value = "email"
$ws.Range("H1").Value = "date_in"
$ws.Range("I1").Value = "date_out"

# --- Sample data row (row 2) ----------------------------------------------
$ws.Range("A2").Value = "C001"
$ws.Range("B2").Value = "B001"

# nik / date_in / date_out look like a number / dates - force text format
# first so the import keeps them as literal text (e.g. doesn't lose leading
# zeros or turn the dates into serial numbers).
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "20240001"

$ws.Range("D2").Value = "mahar"
$ws.Range("E2").Value = "mahatma"
$ws.Range("F2").Value = "mahardhika"
$ws.Range("G2").Value = "test@test.com"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "2024-01-01"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "2024-02-02"
